# Update "想去人数" (interested-people count) values in column F
# for sheets "展览" (Worksheets item 1) and "全部类型" (Worksheets item 4).
# Both sheets share identical data and need identical updates.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F3"  = 2900
    "F5"  = 149
    "F7"  = 1558
    "F11" = 1288
    "F13" = 414
    "F15" = 64
    "F16" = 53
    "F18" = 84
    "F20" = 2920
    "F21" = 356
    "F22" = 14
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
